# Auto update: 2025-12-05 04:00:45
# Updates the hedging/insurance analysis sheet with refreshed model output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$watch = "⛔ 관망하십시오."
$neutral = "⚪ 중립 구간"

# Row 2 - UnitedHealth Group Incorporated (UNH)
$ws.Range("D2").Value = 331.56
$ws.Range("E2").Value = 49.4
$ws.Range("F2").Value = 0.5600000000000001
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 59.4
$ws.Range("M2").Value = $watch
$ws.Range("N2").Value = 54.02451352198364
$ws.Range("O2").Value = $neutral

# Row 3 - MetLife, Inc. (MET)
$ws.Range("F3").Value = 2.41
$ws.Range("K3").Value = 58.2
$ws.Range("M3").Value = $watch
$ws.Range("N3").Value = 54.02451352198364
$ws.Range("O3").Value = $neutral

# Row 4 - American International Group, Inc. (AIG)
$ws.Range("D4").Value = 77.83
$ws.Range("E4").Value = 46.7
$ws.Range("F4").Value = 2.35
$ws.Range("H4").Value = 40
$ws.Range("K4").Value = 50.2
$ws.Range("M4").Value = $watch
$ws.Range("N4").Value = 54.02451352198364
$ws.Range("O4").Value = $neutral

# Row 5 - Prudential Financial, Inc. (PRU)
$ws.Range("D5").Value = 110.76
$ws.Range("E5").Value = 64.90000000000001
$ws.Range("F5").Value = 2.58
$ws.Range("G5").Value = 60
$ws.Range("K5").Value = 50.2
$ws.Range("M5").Value = $watch
$ws.Range("N5").Value = 54.02451352198364
$ws.Range("O5").Value = $neutral

$wb.Save()
